$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.175.14'
$ws.Range('D3').Value = '1.677.97'
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.41'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.519'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.99'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +7.73%  '
$ws.Range('E9').Value = '  +3.28%  '
$ws.Range('E10').Value = '  -0.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0891'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('D12').Value = '1.915.19'
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('D13').Value = '1.681.10'
$ws.Range('E13').Value = '  +0.01%  '
$ws.Range('E14').Value = '  +2.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.559'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.65'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Value = '27.142.21'
$ws.Range('E17').Value = '  +0.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '235.47'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('D19').Value = '0.0₃0743'
$ws.Range('E19').Value = '  +1.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.83'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.78%  '
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.56'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.56'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.45%  '
$ws.Range('E24').Value = '  -0.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.80'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.34%  '
$ws.Range('E26').Value = '  +3.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.44'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('E28').Value = '  +0.41%  '
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('E30').Value = '  +0.54%  '
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.38'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.21%  '
$ws.Range('D33').Value = '1.536.55'
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('E34').Value = '  +1.72%  '
$ws.Range('E35').Value = '  -3.05%  '
$ws.Range('E36').Value = '  +4.43%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.951'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.50%  '
$ws.Range('E38').Value = '  -0.16%  '
$ws.Range('E39').Value = '  -0.78%  '
$ws.Range('E40').Value = '  +2.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '69.91'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.10%  '
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D45').Value = '1.823.32'
$ws.Range('E45').Value = '  +0.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.784'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.54%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.66'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +7.79%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '90.33'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('E49').Value = '  +2.63%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.26'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.64%  '
$ws.Range('E51').Value = '  +1.39%  '
